# Auto-generated edit script: update numeric values in the profit-calculation sheets
# to reflect refreshed market pricing data (per "chore: update Sheets via scheduled runner").
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1293.2273
$ws.Range("I40").Value = 1120.9166
$ws.Range("K40").Value = 1120.9166
$ws.Range("M40").Value = -945.9166
$ws.Range("H45").Value = 2999.5
$ws.Range("J45").Value = 2999.5
$ws.Range("L45").Value = 8998.5
$ws.Range("N45").Value = -9382.5
$ws.Range("H64").Value = 3417.3076
$ws.Range("J64").Value = 3597.8948
$ws.Range("L64").Value = 3597.8948
$ws.Range("N64").Value = -4093.8948
$ws.Range("H67").Value = 3417.3076
$ws.Range("J67").Value = 3597.8948
$ws.Range("L67").Value = 3597.8948
$ws.Range("N67").Value = -5313.8948
$ws.Range("H125").Value = 3186.2
$ws.Range("I125").Value = 3298.3333
$ws.Range("J125").Value = 3018
$ws.Range("K125").Value = 29684.9997
$ws.Range("L125").Value = 27162
$ws.Range("M125").Value = -27224.9997
$ws.Range("N125").Value = -32082
$ws.Range("H129").Value = 893.23914
$ws.Range("J129").Value = 899.75555
$ws.Range("L129").Value = 2699.26665
$ws.Range("N129").Value = -12699.26665
$ws.Range("H132").Value = 4733.4414
$ws.Range("I132").Value = 2712
$ws.Range("J132").Value = 14166.833
$ws.Range("K132").Value = 8136
$ws.Range("L132").Value = 42500.499
$ws.Range("M132").Value = -5606
$ws.Range("N132").Value = -47560.499
$ws.Range("H135").Value = 238.5
$ws.Range("I135").Value = 77.5
$ws.Range("J135").Value = 399.5
$ws.Range("K135").Value = 697.5
$ws.Range("L135").Value = 3595.5
$ws.Range("M135").Value = 1837.5
$ws.Range("N135").Value = -8665.5
$ws.Range("H137").Value = 1212.3429
$ws.Range("I137").Value = 1061.8636
$ws.Range("J137").Value = 1467
$ws.Range("K137").Value = 3185.5908
$ws.Range("L137").Value = 4401
$ws.Range("M137").Value = -635.5907999999999
$ws.Range("N137").Value = -9501
$ws.Range("H141").Value = 10121.583
$ws.Range("I141").Value = 10704.909
$ws.Range("J141").Value = 3705
$ws.Range("K141").Value = 32114.727
$ws.Range("L141").Value = 11115
$ws.Range("M141").Value = -26934.727
$ws.Range("N141").Value = -21475

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1610.7333
$ws.Range("I2").Value = 845.1667
$ws.Range("J2").Value = 4673
$ws.Range("K2").Value = 845.1667
$ws.Range("L2").Value = 4673
$ws.Range("M2").Value = -732.1667
$ws.Range("N2").Value = -4899
$ws.Range("H61").Value = 1736.125
$ws.Range("I61").Value = 1564.8334
$ws.Range("K61").Value = 1564.8334
$ws.Range("M61").Value = -1352.8334
$ws.Range("H74").Value = 662.4318
$ws.Range("I74").Value = 633.17645
$ws.Range("K74").Value = 633.17645
$ws.Range("M74").Value = 240.82355
$ws.Range("H77").Value = 662.4318
$ws.Range("I77").Value = 633.17645
$ws.Range("K77").Value = 3165.88225
$ws.Range("M77").Value = 1202.11775
$ws.Range("H116").Value = 1610.7333
$ws.Range("I116").Value = 845.1667
$ws.Range("J116").Value = 4673
$ws.Range("K116").Value = 845.1667
$ws.Range("L116").Value = 4673
$ws.Range("M116").Value = 1448.8333
$ws.Range("N116").Value = -9261
$ws.Range("H132").Value = 3478.7083
$ws.Range("I132").Value = 3374.55
$ws.Range("K132").Value = 10123.65
$ws.Range("M132").Value = -7593.650000000001
$ws.Range("H136").Value = 1736.125
$ws.Range("I136").Value = 1564.8334
$ws.Range("K136").Value = 4694.5002
$ws.Range("M136").Value = -2144.5002

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1610.7333
$ws.Range("I3").Value = 845.1667
$ws.Range("J3").Value = 4673
$ws.Range("K3").Value = 845.1667
$ws.Range("L3").Value = 4673
$ws.Range("M3").Value = -731.1667
$ws.Range("N3").Value = -4901
$ws.Range("H133").Value = 21920
$ws.Range("J133").Value = 21920
$ws.Range("L133").Value = 21920
$ws.Range("N133").Value = -32040
$ws.Range("H134").Value = 13970.353
$ws.Range("I134").Value = 11245.182
$ws.Range("K134").Value = 33735.546
$ws.Range("M134").Value = -31200.546

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 694.6629
$ws.Range("I31").Value = 617.4
$ws.Range("J31").Value = 854.5172
$ws.Range("K31").Value = 617.4
$ws.Range("L31").Value = 854.5172
$ws.Range("M31").Value = -322.4
$ws.Range("N31").Value = -1444.5172
$ws.Range("H34").Value = 694.6629
$ws.Range("I34").Value = 617.4
$ws.Range("J34").Value = 854.5172
$ws.Range("K34").Value = 617.4
$ws.Range("L34").Value = 854.5172
$ws.Range("M34").Value = -415.4
$ws.Range("N34").Value = -1258.5172
$ws.Range("H42").Value = 10062
$ws.Range("J42").Value = 10062
$ws.Range("L42").Value = 10062
$ws.Range("N42").Value = -11248
$ws.Range("H62").Value = 28574140
$ws.Range("I62").Value = 2494.75
$ws.Range("J62").Value = 66669668
$ws.Range("K62").Value = 2494.75
$ws.Range("L62").Value = 66669668
$ws.Range("M62").Value = -1870.75
$ws.Range("N62").Value = -66670916
$ws.Range("H65").Value = 28574140
$ws.Range("I65").Value = 2494.75
$ws.Range("J65").Value = 66669668
$ws.Range("K65").Value = 12473.75
$ws.Range("L65").Value = 333348340
$ws.Range("M65").Value = -9353.75
$ws.Range("N65").Value = -333354580
$ws.Range("H99").Value = 2121
$ws.Range("J99").Value = 2337.6667
$ws.Range("L99").Value = 2337.6667
$ws.Range("N99").Value = -5333.6667
$ws.Range("H126").Value = 2121
$ws.Range("J126").Value = 2337.6667
$ws.Range("L126").Value = 7013.000100000001
$ws.Range("N126").Value = -11953.0001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1379.3928
$ws.Range("J68").Value = 2218.4614
$ws.Range("L68").Value = 6655.3842
$ws.Range("N68").Value = -8277.3842
$ws.Range("H71").Value = 1379.3928
$ws.Range("J71").Value = 2218.4614
$ws.Range("L71").Value = 19966.1526
$ws.Range("N71").Value = -28078.1526
$ws.Range("H140").Value = 28037.166
$ws.Range("I140").Value = 40791.965
$ws.Range("K140").Value = 122375.895
$ws.Range("M140").Value = -117195.895

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3113.5
$ws.Range("I61").Value = 2977
$ws.Range("J61").Value = 3250
$ws.Range("K61").Value = 2977
$ws.Range("L61").Value = 3250
$ws.Range("M61").Value = -2775
$ws.Range("N61").Value = -3654
$ws.Range("H113").Value = 3113.5
$ws.Range("I113").Value = 2977
$ws.Range("J113").Value = 3250
$ws.Range("K113").Value = 2977
$ws.Range("L113").Value = 3250
$ws.Range("M113").Value = -807
$ws.Range("N113").Value = -7590
$ws.Range("H140").Value = 39914.5
$ws.Range("J140").Value = 39914.5
$ws.Range("L140").Value = 39914.5
$ws.Range("N140").Value = -50274.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 28999.334
$ws.Range("J69").Value = 28999.334
$ws.Range("L69").Value = 28999.334
$ws.Range("N69").Value = -30497.334
$ws.Range("H72").Value = 28999.334
$ws.Range("J72").Value = 28999.334
$ws.Range("L72").Value = 86998.00199999999
$ws.Range("N72").Value = -94486.00199999999
$ws.Range("H136").Value = 2088.875
$ws.Range("I136").Value = 1934.3334
$ws.Range("K136").Value = 5803.0002
$ws.Range("M136").Value = -3253.0002
